$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Average" column (O) and "Active" column (N) are being added after the
# existing M ("Mother's Name") column. Pick up the header formatting (bold
# font + fill) from the neighboring header cell (M1) before writing the new
# header text, so the new headers look the same as the rest of row 1.
# Note: O1 ("Average") is written first so it lands on shared-string index
# 35, followed by N1 ("Active") on index 36.
$ws.Range("M1").Copy()
$ws.Range("O1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("O1").Value = "Average"

$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("N1").Value = "Active"
$excel.CutCopyMode = $false

# "Active" flag per student (boolean column N)
$ws.Range("N2").Value = $true
$ws.Range("N3").Value = $false
$ws.Range("N4").Value = $true
$ws.Range("N5").Value = $true

# "Average" score per student (numeric column O, 2 decimal places)
$ws.Range("O2").Value = 23.41
$ws.Range("O3").Value = 21.67
$ws.Range("O4").Value = 29.32
$ws.Range("O5").Value = 22.98
$ws.Range("O2:O5").NumberFormat = "0.00"

# Match the bestFit-computed column widths for the two new columns.
$ws.Columns.Item(14).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 6.998697916666667

# Leave the same selection state captured in the saved workbook.
[void]$ws.Range("N2").Select()
